# Expense Tracker Project - completed
# - Rename sheet "expense" -> "Expense"
# - Populate the expense table (Icon, Category, Amount, Date)
# - Format the Date column as a short date (numFmtId 14)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet title fix
$ws.Name = "Expense"

# Header row
$ws.Range("A1").Value = "Icon"
$ws.Range("B1").Value = "Category"
$ws.Range("C1").Value = "Amount"
$ws.Range("D1").Value = "Date"

# Row data: Icon, Category, Amount, Date (OLE Automation date serial)
$rows = @(
    @("💡", "Utilities",         1000,  45862.22928240741),
    @("🛒", "me",                1200,  45862.22928240741),
    @("🛒", "Groceries",         2500,  45861.22928240741),
    @("💡", "Electricity Bill", 13000,  45845.22928240741),
    @("🎬", "Entertainment",      750,  45843.22928240741),
    @("🏠", "Rent",              9500,  45841.22928240741),
    @("🚗", "Transport",         1200,  45840.22928240741),
    @("🚕", "Transport",          120,  45840.22928240741),
    @("🍔", "Food",               450,  45839.22928240741),
    @("🍕", "Food",               350,  45839.22928240741)
)

$firstRow = 2
$lastRow = $firstRow + $rows.Count - 1

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $firstRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}

# Set the first date cell's value + number format, then propagate the SAME
# style to the rest of the column via copy/paste-format (avoids minting a
# brand new style record per cell) before filling in their values.
$ws.Cells.Item($firstRow, 4).Value = $rows[0][3]
$ws.Cells.Item($firstRow, 4).NumberFormat = "mm-dd-yy"

$ws.Cells.Item($firstRow, 4).Copy()
$ws.Range($ws.Cells.Item($firstRow + 1, 4), $ws.Cells.Item($lastRow, 4)).PasteSpecial(-4122)

for ($i = 1; $i -lt $rows.Count; $i++) {
    $ws.Cells.Item($firstRow + $i, 4).Value = $rows[$i][3]
}
